{"js": "// The student/group identifier on the title page was corrected from\n// \"3M\" to \"7N\". In the underlying OOXML this lives as two separate\n// runs (\"3\" then \"M\"); after the edit they collapse into a single run\n// containing \"7N\" while the earlier run (\"22310439           \") is\n// left untouched.\nconst body = context.document.body;\n\nconst results = body.search(\"3M\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\"7N\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# The student/group identifier on the title page was corrected from\n# \"3M\" to \"7N\" (e.g. \"22310439           3M\" -> \"22310439           7N\").\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Text = \"3M\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.Replacement.Text = \"7N\"\n$find.Execute([ref]$find.Text, $true, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n"}
